$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'27.635.66"
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.Value = "'  +0.09%  "
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.Value = "'1.634.44"
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.Value = "'  -0.11%  "
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.Value = "'  -0.13%  "
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.Value = "'212.29"
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.Value = "'  -0.02%  "
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.Value = "'0.521"
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.Value = "'  -0.48%  "
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.Value = "'  -0.15%  "
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.Value = "'23.33"
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.Value = "'  +1.65%  "
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.Value = "'  +2.60%  "
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.Value = "'  +0.33%  "
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.Value = "'0.0872"
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.Value = "'  -2.34%  "
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.Value = "'1.864.11"
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.Value = "'  -0.22%  "
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.Value = "'1.639.86"
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.Value = "'  +0.26%  "
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.Value = "'  +0.43%  "
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.Value = "'0.554"
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.Value = "'  -0.81%  "
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.Value = "'65.36"
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.Value = "'  +1.29%  "
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.Value = "'27.606.85"
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.Value = "'  +0.00%  "
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.Value = "'231.65"
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.Value = "'  +1.27%  "
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.Value = "'0.0₃0721"
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.Value = "'  -0.10%  "
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.Value = "'7.59"
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.Value = "'  -1.53%  "
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.Value = "'1.00"
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.Value = "'10.65"
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.Value = "'  +6.45%  "
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.Value = "'4.36"
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.Value = "'  +1.71%  "
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.Value = "'  +8.03%  "
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.Value = "'149.72"
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.Value = "'  -0.64%  "
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.Value = "'  -0.62%  "
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.Value = "'  -0.13%  "
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.Value = "'15.56"
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.Value = "'  -0.13%  "
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.Value = "'  -0.10%  "
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.Value = "'  -0.34%  "
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.Value = "'0.0485"
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.Value = "'  +0.01%  "
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.Value = "'  -0.40%  "
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.Value = "'1.480.09"
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.Value = "'  +1.73%  "
$c.Style = 'Normal'
$c = $ws.Range('D34')
$c.Value = "'3.08"
$c.Style = 'Normal'
$c = $ws.Range('E34')
$c.Value = "'  -0.79%  "
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.Value = "'  -1.01%  "
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.Value = "'  -1.37%  "
$c.Style = 'Normal'
$c = $ws.Range('E37')
$c.Value = "'  +5.22%  "
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.Value = "'0.882"
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.Value = "'  +0.59%  "
$c.Style = 'Normal'
$c = $ws.Range('B39')
$c.Value = "'ImmutableX"
$c.Style = 'Normal'
$c = $ws.Range('C39')
$c.Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.Value = "'0.560"
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.Value = "'  -0.83%  "
$c.Style = 'Normal'
$c = $ws.Range('B40')
$c.Value = "'VeChain"
$c.Style = 'Normal'
$c = $ws.Range('C40')
$c.Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.Value = "'0.0167"
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.Value = "'  +0.51%  "
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.Value = "'1.03"
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.Value = "'  +2.33%  "
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.Value = "'  -0.10%  "
$c.Style = 'Normal'
$c = $ws.Range('D43')
$c.Value = "'68.22"
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.Value = "'  -2.04%  "
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.Value = "'  -0.09%  "
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.Value = "'  -0.79%  "
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.Value = "'  -4.31%  "
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.Value = "'1.773.98"
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.Value = "'  -0.27%  "
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.Value = "'1.76"
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.Value = "'  +2.15%  "
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.Value = "'87.72"
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.Value = "'  +1.64%  "
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.Value = "'0.0₆0105"
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.Value = "'  -1.50%  "
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.Value = "'  +1.08%  "
$c.Style = 'Normal'
